$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for columns L, M, N
$ws.Range("L1").Value = "apoio_medio"
$ws.Range("M1").Value = "contribuicoes"
$ws.Range("N1").Value = "media_contribuicoes"

# New data values for columns L, M, N (rows 2-16)
$lValues = @(94.08517009767365, 129.5537663642677, 79.80942858649695, 87.19342470373856, 50.12134015913439, 95.44644410600942, 96.51058441972074, 66.65034280439198, 86.43541554443971, 55.58374799260083, 30.9975069667077, 40.89129143626957, 18.08263434560471, 15.74416694302886, 20.96281755102498)
$mValues = @(7547, 32860, 48629, 174471, 46, 15501, 95943, 17194, 74806, 202, 37, 25, 307, 591, 1248)
$nValues = @(260.2413793103448, 395.9036144578313, 347.35, 302.9010416666667, 23, 224.6521739130435, 218.0522727272727, 97.69318181818181, 108.2575976845152, 28.85714285714286, 5.285714285714286, 12.5, 17.05555555555556, 23.64, 12.48)

for ($i = 0; $i -lt 15; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 12).Value = $lValues[$i]
    $ws.Cells.Item($row, 13).Value = $mValues[$i]
    $ws.Cells.Item($row, 14).Value = $nValues[$i]
}

# Multiply existing E and F columns (taxa_sucesso percentages) by 100
for ($row = 2; $row -le 16; $row++) {
    $eCell = $ws.Cells.Item($row, 5)
    $eCell.Value = $eCell.Value() * 100
    $fCell = $ws.Cells.Item($row, 6)
    $fCell.Value = $fCell.Value() * 100
}
